$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 183.27272
$ws.Range("I55").Value = 148
$ws.Range("J55").Value = 212.66667
$ws.Range("K55").Value = 148
$ws.Range("L55").Value = 212.66667
$ws.Range("M55").Value = 66
$ws.Range("N55").Value = -640.6666700000001
$ws.Range("H138").Value = 3877.3738
$ws.Range("I138").Value = 1527.7587
$ws.Range("J138").Value = 4850.7856
$ws.Range("K138").Value = 4583.2761
$ws.Range("L138").Value = 14552.3568
$ws.Range("M138").Value = 556.7239
$ws.Range("N138").Value = -24832.3568

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 500
$ws.Range("I26").Value = 500
$ws.Range("K26").Value = 500
$ws.Range("M26").Value = -170
$ws.Range("H44").Value = 26814
$ws.Range("J44").Value = 35139.6
$ws.Range("L44").Value = 35139.6
$ws.Range("N44").Value = -36115.6
$ws.Range("H55").Value = 34946.332
$ws.Range("J55").Value = 34946.332
$ws.Range("L55").Value = 34946.332
$ws.Range("N55").Value = -35576.332
$ws.Range("H80").Value = 37769.332
$ws.Range("J80").Value = 37769.332
$ws.Range("L80").Value = 37769.332
$ws.Range("N80").Value = -39765.332
$ws.Range("H83").Value = 37769.332
$ws.Range("J83").Value = 37769.332
$ws.Range("L83").Value = 113307.996
$ws.Range("N83").Value = -123291.996
$ws.Range("H88").Value = 2385.818
$ws.Range("I88").Value = 3099.1
$ws.Range("J88").Value = 1791.4166
$ws.Range("K88").Value = 3099.1
$ws.Range("L88").Value = 1791.4166
$ws.Range("M88").Value = -2693.1
$ws.Range("N88").Value = -2603.4166
$ws.Range("H91").Value = 2385.818
$ws.Range("I91").Value = 3099.1
$ws.Range("J91").Value = 1791.4166
$ws.Range("K91").Value = 3099.1
$ws.Range("L91").Value = 1791.4166
$ws.Range("M91").Value = -1695.1
$ws.Range("N91").Value = -4599.4166
$ws.Range("H122").Value = 14131.823
$ws.Range("I122").Value = 22104.2
$ws.Range("K122").Value = 66312.60000000001
$ws.Range("M122").Value = -63862.60000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35372
$ws.Range("J35").Value = 35372
$ws.Range("L35").Value = 35372
$ws.Range("N35").Value = -35992
$ws.Range("H82").Value = 20453.354
$ws.Range("J82").Value = 37967.75
$ws.Range("L82").Value = 37967.75
$ws.Range("N82").Value = -38733.75
$ws.Range("H85").Value = 20453.354
$ws.Range("J85").Value = 37967.75
$ws.Range("L85").Value = 37967.75
$ws.Range("N85").Value = -40619.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3269764
$ws.Range("I16").Value = 5883375
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 5883375
$ws.Range("L16").Value = 2750
$ws.Range("M16").Value = -5883088
$ws.Range("N16").Value = -3324
$ws.Range("H18").Value = 37000
$ws.Range("J18").Value = 37000
$ws.Range("L18").Value = 37000
$ws.Range("N18").Value = -37460
$ws.Range("H41").Value = 20598.5
$ws.Range("J41").Value = 28946
$ws.Range("L41").Value = 28946
$ws.Range("N41").Value = -29802
$ws.Range("H50").Value = 8916
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 8916
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 8916
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -10166
$ws.Range("H51").Value = 9371.333000000001
$ws.Range("J51").Value = 9371.333000000001
$ws.Range("L51").Value = 9371.333000000001
$ws.Range("N51").Value = -10843.333
$ws.Range("H61").Value = 9371.333000000001
$ws.Range("J61").Value = 9371.333000000001
$ws.Range("L61").Value = 9371.333000000001
$ws.Range("N61").Value = -10067.333
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 16907.5
$ws.Range("I68").Value = 15000
$ws.Range("J68").Value = 17289
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 17289
$ws.Range("M68").Value = -14251
$ws.Range("N68").Value = -18787
$ws.Range("H71").Value = 16907.5
$ws.Range("I71").Value = 15000
$ws.Range("J71").Value = 17289
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 51867
$ws.Range("M71").Value = -41256
$ws.Range("N71").Value = -59355
$ws.Range("H86").Value = 41672800
$ws.Range("I86").Value = 66674764
$ws.Range("J86").Value = 2857.111
$ws.Range("K86").Value = 66674764
$ws.Range("L86").Value = 2857.111
$ws.Range("M86").Value = -66673641
$ws.Range("N86").Value = -5103.111
$ws.Range("H89").Value = 41672800
$ws.Range("I89").Value = 66674764
$ws.Range("J89").Value = 2857.111
$ws.Range("K89").Value = 333373820
$ws.Range("L89").Value = 14285.555
$ws.Range("M89").Value = -333368204
$ws.Range("N89").Value = -25517.555
$ws.Range("H109").Value = 13282.25
$ws.Range("J109").Value = 13282.25
$ws.Range("L109").Value = 13282.25
$ws.Range("N109").Value = -15362.25
$ws.Range("H113").Value = 3269764
$ws.Range("I113").Value = 5883375
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 5883375
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = -5881205
$ws.Range("N113").Value = -7090
$ws.Range("H120").Value = 12939.637
$ws.Range("J120").Value = 12939.637
$ws.Range("L120").Value = 12939.637
$ws.Range("N120").Value = -20197.637
$ws.Range("H122").Value = 1984
$ws.Range("I122").Value = 1376
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 4128
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -1678
$ws.Range("N122").Value = -14500

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 7350
$ws.Range("I110").Value = 6300
$ws.Range("J110").Value = 7980
$ws.Range("K110").Value = 18900
$ws.Range("L110").Value = 23940
$ws.Range("M110").Value = -14810
$ws.Range("N110").Value = -32120

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14963.059
$ws.Range("I57").Value = 12990
$ws.Range("J57").Value = 15086.375
$ws.Range("K57").Value = 12990
$ws.Range("L57").Value = 15086.375
$ws.Range("M57").Value = -12170
$ws.Range("N57").Value = -16726.375
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
$ws.Range("H70").Value = 30610.125
$ws.Range("I70").Value = 38190.582
$ws.Range("J70").Value = 4499.6665
$ws.Range("K70").Value = 38190.582
$ws.Range("L70").Value = 4499.6665
$ws.Range("M70").Value = -37920.582
$ws.Range("N70").Value = -5039.6665
$ws.Range("H73").Value = 30610.125
$ws.Range("I73").Value = 38190.582
$ws.Range("J73").Value = 4499.6665
$ws.Range("K73").Value = 38190.582
$ws.Range("L73").Value = 4499.6665
$ws.Range("M73").Value = -37254.582
$ws.Range("N73").Value = -6371.6665
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H82").Value = 40666.668
$ws.Range("J82").Value = 40666.668
$ws.Range("L82").Value = 40666.668
$ws.Range("N82").Value = -41432.668
$ws.Range("H85").Value = 40666.668
$ws.Range("J85").Value = 40666.668
$ws.Range("L85").Value = 40666.668
$ws.Range("N85").Value = -43318.668
$ws.Range("H122").Value = 2489.5264
$ws.Range("I122").Value = 2391
$ws.Range("J122").Value = 2625
$ws.Range("K122").Value = 7173
$ws.Range("L122").Value = 7875
$ws.Range("M122").Value = -4723
$ws.Range("N122").Value = -12775
$ws.Range("H126").Value = 6373.5
$ws.Range("I126").Value = 2830
$ws.Range("J126").Value = 7306
$ws.Range("K126").Value = 8490
$ws.Range("L126").Value = 21918
$ws.Range("M126").Value = -6020
$ws.Range("N126").Value = -26858

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4449.4346
$ws.Range("I122").Value = 3837.0908
$ws.Range("J122").Value = 5010.75
$ws.Range("K122").Value = 11511.2724
$ws.Range("L122").Value = 15032.25
$ws.Range("M122").Value = -9061.2724
$ws.Range("N122").Value = -19932.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1675.8966
$ws.Range("I107").Value = 1272.3158
$ws.Range("K107").Value = 3816.9474
$ws.Range("M107").Value = -1896.9474
$ws.Range("H109").Value = 21000.334
$ws.Range("J109").Value = 21000.334
$ws.Range("L109").Value = 21000.334
$ws.Range("N109").Value = -23774.334
$ws.Range("H122").Value = 14634.588
$ws.Range("I122").Value = 23489.777
$ws.Range("J122").Value = 4672.5
$ws.Range("K122").Value = 70469.33099999999
$ws.Range("L122").Value = 14017.5
$ws.Range("M122").Value = -68019.33099999999
$ws.Range("N122").Value = -18917.5
